# Refresh the scraped coin price/volume figures (GitHub Actions cron update).
# Price (col D) and Volume(1h) (col E) are stored as plain text in the sheet.
# For D-values that look like plain numbers ("692.95", "1.00", "0.180", ...)
# a leading apostrophe forces Excel to keep them as literal text (preserving
# exact digits/trailing zeros instead of silently re-parsing them as floats);
# the style is then reset to "Normal" so no stray quote-prefix formatting is
# left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.240.77'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '3.866.82'
$ws.Range("E3").Value = '  +1.52%  '
$ws.Range("D5").Value = '''692.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.06%  '
$ws.Range("D6").Value = '''173.26'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.24%  '
$ws.Range("D7").Value = '3.865.05'
$ws.Range("E7").Value = '  +1.49%  '
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("E10").Value = '  +1.68%  '
$ws.Range("D11").Value = '''7.38'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.03%  '
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("E13").Value = '  +6.22%  '
$ws.Range("D14").Value = '''36.67'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.40%  '
$ws.Range("D15").Value = '4.516.93'
$ws.Range("E15").Value = '  +1.52%  '
$ws.Range("D16").Value = '3.868.19'
$ws.Range("E16").Value = '  +1.62%  '
$ws.Range("D17").Value = '71.273.42'
$ws.Range("E17").Value = '  +1.24%  '
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("D19").Value = '''7.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.05%  '
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("D21").Value = '''11.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("D22").Value = '''494.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.56%  '
$ws.Range("E23").Value = '  +1.45%  '
$ws.Range("D24").Value = '''84.95'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.69%  '
$ws.Range("E25").Value = '  +3.58%  '
$ws.Range("E26").Value = '  +1.34%  '
$ws.Range("D27").Value = '''10.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.08%  '
$ws.Range("E28").Value = '  +1.83%  '
$ws.Range("D29").Value = '4.018.36'
$ws.Range("E29").Value = '  +1.51%  '
$ws.Range("D30").Value = '''1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").Value = '''3.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +9.83%  '
$ws.Range("D32").Value = '''7.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.61%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("D34").Value = '''29.80'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.51%  '
$ws.Range("D35").Value = '''0.180'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.81%  '
$ws.Range("E36").Value = '  +2.25%  '
$ws.Range("D37").Value = '3.817.43'
$ws.Range("E37").Value = '  +1.41%  '
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("E39").Value = '  +2.73%  '
$ws.Range("E40").Value = '  +13.73%  '
$ws.Range("E41").Value = '  +1.16%  '
$ws.Range("D42").Value = '''6.08'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.91%  '
$ws.Range("D43").Value = '''1.03'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.07%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D46").Value = '''164.06'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.34%  '
$ws.Range("D47").Value = '''0.000309'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.34%  '
$ws.Range("E48").Value = '  +1.04%  '
$ws.Range("D49").Value = '''44.47'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.20%  '
$ws.Range("D50").Value = '''0.304'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.09%  '
$ws.Range("D51").Value = '''8.71'
$ws.Range("D51").Style = "Normal"
